$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.971.13'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.635.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.37%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.06'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.524'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.40'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.97%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.98%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0614'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.866.96'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.640.26'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.563'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.77'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.962.87'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.11'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.65'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.52%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.90%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.56%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.37'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.09%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.16%  '

$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.111'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.52%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.65'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.57%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.33%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.93%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.02%  '

$ws.Range("B33").Value = 'Maker'
$ws.Range("C33").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.403.49'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.33%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.07'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +12.67%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.73%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.86%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.865'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.94%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.66%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.81'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.83'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.56%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.776.97'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.27'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.73%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.94%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.27%  '
